# Update "想去人数" (interested-count) column F on the 展览 (sheet1) and
# 全部类型 (sheet4) sheets, which carry duplicate rows for the same events.
# Values below reflect a fresh scrape run (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# row -> (new value on 展览/全部类型 sheets other than row 12,
#         row 12 differs slightly between the two sheets)
$changes = @{
    2  = 361
    4  = 10714
    5  = 330
    6  = 972
    7  = 151
    8  = 1324
    9  = 8264
    10 = 35
    13 = 217
    15 = 3291
    18 = 774
    19 = 129
    20 = 1064
    23 = 1754
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $changes.Keys) {
        $ws.Cells.Item($row, 6).Value = $changes[$row]
    }
}

# Row 12 ("F12") gets a different value on each of the two sheets.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(12, 6).Value = 168

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(12, 6).Value = 169
